# Automatic update of files: refresh the "Förändrad" (Changed) date column
# for rows 2-16 from serial 45221 (2023-10-22) to serial 45224 (2023-10-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45221) {
        $cell.Value2 = 45224
    }
}
